# Refresh the cryptocurrency price/volume snapshot with the latest values
# pulled in by the scheduled GitHub Actions scraping workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.569.43'
$ws.Range("E2").Value = '  -0.85%  '

# Row 3
$ws.Range("D3").Value = '2.512.96'
$ws.Range("E3").Value = '  -1.47%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$c = $ws.Range("D5")
$c.Value = "'317.25"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +4.34%  '

# Row 6
$c = $ws.Range("D6")
$c.Value = "'95.24"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -3.12%  '

# Row 7
$ws.Range("E7").Value = '  +0.38%  '

# Row 8
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("E9").Value = '  -1.92%  '

# Row 10
$c = $ws.Range("D10")
$c.Value = "'35.98"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.78%  '

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.0807"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.25%  '

# Row 12
$c = $ws.Range("D12")
$c.Value = "'7.68"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.76%  '

# Row 14
$ws.Range("D14").Value = '2.896.22'
$ws.Range("E14").Value = '  -1.55%  '

# Row 15
$c = $ws.Range("D15")
$c.Value = "'15.47"
$c.Style = "Normal"
$ws.Range("E15").Value = '  +4.49%  '

# Row 16
$ws.Range("D16").Value = '2.530.00'
$ws.Range("E16").Value = '  -2.44%  '

# Row 17
$c = $ws.Range("D17")
$c.Value = "'0.853"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.82%  '

# Row 18
$ws.Range("D18").Value = '42.591.93'
$ws.Range("E18").Value = '  -1.28%  '

# Row 19
$c = $ws.Range("D19")
$c.Value = "'13.10"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -3.99%  '

# Row 20
$ws.Range("E20").Value = '  -1.94%  '

# Row 21
$ws.Range("E21").Value = '  -1.20%  '

# Row 22
$c = $ws.Range("D22")
$c.Value = "'71.20"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.95%  '

# Row 23
$c = $ws.Range("D23")
$c.Value = "'250.68"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.74%  '

# Row 24
$ws.Range("E24").Value = '  +0.90%  '

# Row 25
$c = $ws.Range("D25")
$c.Value = "'2.02"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.11%  '

# Row 26
$c = $ws.Range("D26")
$c.Value = "'26.69"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -4.40%  '

# Row 27
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$c = $ws.Range("D28")
$c.Value = "'2.37"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +12.78%  '

# Row 29
$c = $ws.Range("D29")
$c.Value = "'38.69"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.31%  '

# Row 30
$c = $ws.Range("D30")
$c.Value = "'10.03"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.01%  '

# Row 31
$c = $ws.Range("D31")
$c.Value = "'5.89"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.46%  '

# Row 32
$c = $ws.Range("D32")
$c.Value = "'156.11"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.97%  '

# Row 33
$c = $ws.Range("D33")
$c.Value = "'19.55"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +3.52%  '

# Row 34
$ws.Range("E34").Value = '  +0.47%  '

# Row 35
$ws.Range("E35").Value = '  -3.80%  '

# Row 36
$c = $ws.Range("D36")
$c.Value = "'0.0785"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.43%  '

# Row 37
$ws.Range("E37").Value = '  -5.33%  '

# Row 38
$ws.Range("E38").Value = '  -1.72%  '

# Row 39
$ws.Range("E39").Value = '  -0.09%  '

# Row 40
$c = $ws.Range("D40")
$c.Value = "'23.95"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -7.02%  '

# Row 41
$ws.Range("E41").Value = '  +0.58%  '

# Row 42
$c = $ws.Range("D42")
$c.Value = "'3.85"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.85%  '

# Row 43
$c = $ws.Range("D43")
$c.Value = "'3.37"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.83%  '

# Row 44
$ws.Range("E44").Value = '  +0.08%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.046.80'
$ws.Range("E45").Value = '  -2.54%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D46")
$c.Value = "'0.0299"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.12%  '

# Row 47
$c = $ws.Range("D47")
$c.Value = "'84.33"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -2.27%  '

# Row 48
$c = $ws.Range("D48")
$c.Value = "'8.79"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.38%  '

# Row 49
$ws.Range("D49").Value = '2.753.51'
$ws.Range("E49").Value = '  -1.65%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D50")
$c.Value = "'0.189"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.97%  '

# Row 51
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$c = $ws.Range("D51")
$c.Value = "'72.43"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.32%  '
